$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.019738287951536
$ws.Range("D2").Value = 1.025030341685802
$ws.Range("E2").Value = 1.020836431697836
$ws.Range("F2").Value = 1.030944436962417
$ws.Range("I2").Value = 1.029496139713008
$ws.Range("J2").Value = 1.02493935506123
$ws.Range("K2").Value = 1.027857106925973
$ws.Range("L2").Value = 1.023675539495702
$ws.Range("M2").Value = 1.033753982721449
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.020703686673239
$ws.Range("D3").Value = 1.025726643926148
$ws.Range("E3").Value = 1.021655332347686
$ws.Range("F3").Value = 1.032136402432164
$ws.Range("I3").Value = 1.029689213294664
$ws.Range("J3").Value = 1.025541555563746
$ws.Range("K3").Value = 1.028360927276579
$ws.Range("L3").Value = 1.024300730009215
$ws.Range("M3").Value = 1.034753379959492
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.021328463415899
$ws.Range("D4").Value = 1.026176907473451
$ws.Range("E4").Value = 1.022185701889502
$ws.Range("F4").Value = 1.032907738130818
$ws.Range("I4").Value = 1.029812313717949
$ws.Range("J4").Value = 1.025930745148472
$ws.Range("K4").Value = 1.028685970249506
$ws.Range("L4").Value = 1.024705110954109
$ws.Range("M4").Value = 1.035399548077043
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.021591142922445
$ws.Range("D5").Value = 1.026366127659805
$ws.Range("E5").Value = 1.022408784695844
$ws.Range("F5").Value = 1.033232021099601
$ws.Range("I5").Value = 1.029863626293857
$ws.Range("J5").Value = 1.026094246310487
$ws.Range("K5").Value = 1.028822387181625
$ws.Range("L5").Value = 1.024875073954772
$ws.Range("M5").Value = 1.035671075563872
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.021635249326208
$ws.Range("D6").Value = 1.026397894390559
$ws.Range("E6").Value = 1.022446248024102
$ws.Range("F6").Value = 1.033286470437726
$ws.Range("I6").Value = 1.029872216153353
$ws.Range("J6").Value = 1.026121692165338
$ws.Range("K6").Value = 1.028845278610151
$ws.Range("L6").Value = 1.024903609199129
$ws.Range("M6").Value = 1.035716659099356
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.021331973261006
$ws.Range("D7").Value = 1.026179436121619
$ws.Range("E7").Value = 1.022188682280392
$ws.Range("F7").Value = 1.032912071162373
$ws.Range("I7").Value = 1.029813001084054
$ws.Range("J7").Value = 1.02593293030821
$ws.Range("K7").Value = 1.028687793968336
$ws.Range("L7").Value = 1.024707382161129
$ws.Range("M7").Value = 1.035403176718142
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020064528193567
$ws.Range("D8").Value = 1.025265719823458
$ws.Range("E8").Value = 1.021113081816247
$ws.Range("F8").Value = 1.03134725628497
$ws.Range("I8").Value = 1.029561768690441
$ws.Range("J8").Value = 1.025142969431612
$ws.Range("K8").Value = 1.028027574234045
$ws.Range("L8").Value = 1.023886858112236
$ws.Range("M8").Value = 1.034091839153248
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.017831904535958
$ws.Range("D9").Value = 1.02365345532043
$ws.Range("E9").Value = 1.019221498950827
$ws.Range("F9").Value = 1.028590248579319
$ws.Range("I9").Value = 1.029105061681927
$ws.Range("J9").Value = 1.023747348534171
$ws.Range("K9").Value = 1.026856842232393
$ws.Range("L9").Value = 1.02243979969793
$ws.Range("M9").Value = 1.031777188906146
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.016344030523937
$ws.Range("D10").Value = 1.022577206797014
$ws.Range("E10").Value = 1.01796302558609
$ws.Range("F10").Value = 1.026752478039145
$ws.Range("I10").Value = 1.02879120178292
$ws.Range("J10").Value = 1.022814539440733
$ws.Range("K10").Value = 1.026071459071635
$ws.Range("L10").Value = 1.021474330670373
$ws.Range("M10").Value = 1.030231445588204
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.015699896425218
$ws.Range("D11").Value = 1.022110858042479
$ws.Range("E11").Value = 1.017418715576875
$ws.Range("F11").Value = 1.025956750587812
$ws.Range("I11").Value = 1.028653075352977
$ws.Range("J11").Value = 1.022410060726851
$ws.Range("K11").Value = 1.025730226370975
$ws.Range("L11").Value = 1.021056098102066
$ws.Range("M11").Value = 1.029561489867695
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01546065530643
$ws.Range("D12").Value = 1.021937587281705
$ws.Range("E12").Value = 1.01721662812849
$ws.Range("F12").Value = 1.025661187059344
$ws.Range("I12").Value = 1.028601435420978
$ws.Range("J12").Value = 1.022259734669951
$ws.Range("K12").Value = 1.025603304157358
$ws.Range("L12").Value = 1.020900721670035
$ws.Range("M12").Value = 1.029312541983181
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.015511972477741
$ws.Range("D13").Value = 1.021974756608696
$ws.Range("E13").Value = 1.017259972327397
$ws.Range("F13").Value = 1.025724586223156
$ws.Range("I13").Value = 1.028612527453252
$ws.Range("J13").Value = 1.022291983943561
$ws.Range("K13").Value = 1.025630537231348
$ws.Range("L13").Value = 1.020934051620927
$ws.Range("M13").Value = 1.029365946525158
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.015680120291893
$ws.Range("D14").Value = 1.022096536402676
$ws.Range("E14").Value = 1.017402009052943
$ws.Range("F14").Value = 1.025932319115601
$ws.Range("I14").Value = 1.028648813583898
$ws.Range("J14").Value = 1.022397636452607
$ws.Range("K14").Value = 1.025719738471046
$ws.Range("L14").Value = 1.02104325516988
$ws.Range("M14").Value = 1.029540913729659
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.015783724227672
$ws.Range("D15").Value = 1.02217156261059
$ws.Range("E15").Value = 1.01748953498373
$ws.Range("F15").Value = 1.0260603108545
$ws.Range("I15").Value = 1.028671126468565
$ws.Range("J15").Value = 1.022462721245388
$ws.Range("K15").Value = 1.025774675365808
$ws.Range("L15").Value = 1.021110535619872
$ws.Range("M15").Value = 1.029648703995334
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016386782233533
$ws.Range("D16").Value = 1.02260815003692
$ws.Range("E16").Value = 1.017999162721121
$ws.Range("F16").Value = 1.026805288607026
$ws.Range("I16").Value = 1.028800321973304
$ws.Range("J16").Value = 1.022841371466759
$ws.Range("K16").Value = 1.026094081242707
$ws.Range("L16").Value = 1.021502083680833
$ws.Range("M16").Value = 1.030275894840318
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.016765098212869
$ws.Range("D17").Value = 1.022881923152811
$ws.Range("E17").Value = 1.018319004553044
$ws.Range("F17").Value = 1.027272603411322
$ws.Range("I17").Value = 1.028880767988118
$ws.Range("J17").Value = 1.023078737447903
$ws.Range("K17").Value = 1.026294126859509
$ws.Range("L17").Value = 1.021747644121872
$ws.Range("M17").Value = 1.030669143713461
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.016985775589159
$ws.Range("D18").Value = 1.02304157886023
$ws.Range("E18").Value = 1.018505622386642
$ws.Range("F18").Value = 1.027545184080056
$ws.Range("I18").Value = 1.028927476274301
$ws.Range("J18").Value = 1.023217134349091
$ws.Range("K18").Value = 1.026410698449219
$ws.Range("L18").Value = 1.021890858152786
$ws.Range("M18").Value = 1.030898457543139
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01706102291337
$ws.Range("D19").Value = 1.023096011956427
$ws.Range("E19").Value = 1.018569264314186
$ws.Range("F19").Value = 1.027638127722175
$ws.Range("I19").Value = 1.028943366196717
$ws.Range("J19").Value = 1.023264314831975
$ws.Range("K19").Value = 1.026450427396257
$ws.Range("L19").Value = 1.021939687499553
$ws.Range("M19").Value = 1.030976637203832
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0167245072393
$ws.Range("D20").Value = 1.022852553141931
$ws.Range("E20").Value = 1.018284682390856
$ws.Range("F20").Value = 1.027222464563033
$ws.Range("I20").Value = 1.028872159077529
$ws.Range("J20").Value = 1.023053275975167
$ws.Range("K20").Value = 1.026272675400665
$ws.Range("L20").Value = 1.021721299587774
$ws.Range("M20").Value = 1.030626958204817
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.015630604426773
$ws.Range("D21").Value = 1.022060676608264
$ws.Range("E21").Value = 1.017360180189406
$ws.Range("F21").Value = 1.025871146822926
$ws.Range("I21").Value = 1.02863813742968
$ws.Range("J21").Value = 1.022366526755538
$ws.Range("K21").Value = 1.025693475705704
$ws.Range("L21").Value = 1.021011098166544
$ws.Range("M21").Value = 1.029489392917172
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.014942934119008
$ws.Range("D22").Value = 1.021562514736038
$ws.Range("E22").Value = 1.016779450451999
$ws.Range("F22").Value = 1.025021548684837
$ws.Range("I22").Value = 1.028489068649543
$ws.Range("J22").Value = 1.021934250689791
$ws.Range("K22").Value = 1.025328308094179
$ws.Range("L22").Value = 1.020564414818026
$ws.Range("M22").Value = 1.028773601876947
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.015307470698635
$ws.Range("D23").Value = 1.021826625786133
$ws.Range("E23").Value = 1.01708725468484
$ws.Range("F23").Value = 1.025471934266202
$ws.Range("I23").Value = 1.02856827566133
$ws.Range("J23").Value = 1.02216345467964
$ws.Range("K23").Value = 1.02552198512141
$ws.Range("L23").Value = 1.020801224301613
$ws.Range("M23").Value = 1.02915310938299
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.016742848515646
$ws.Range("D24").Value = 1.022865824283058
$ws.Range("E24").Value = 1.018300190914734
$ws.Range("F24").Value = 1.027245120138137
$ws.Range("I24").Value = 1.028876049736453
$ws.Range("J24").Value = 1.023064781087575
$ws.Range("K24").Value = 1.026282368736473
$ws.Range("L24").Value = 1.021733203601942
$ws.Range("M24").Value = 1.030646020210326
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01840899721076
$ws.Range("D25").Value = 1.02407051643096
$ws.Range("E25").Value = 1.019710067185877
$ws.Range("F25").Value = 1.029302957469763
$ws.Range("I25").Value = 1.029224788115333
$ws.Range("J25").Value = 1.024108574423847
$ws.Range("K25").Value = 1.027160370238551
$ws.Range("L25").Value = 1.022814036219345
$ws.Range("M25").Value = 1.032376046083951
